# "Blank database template.xlsx" update:
#   Variants_SNP_indel sheet gains a new "Reference/s" column (H) with a
#   header comment, plus the saved selection moves to H18.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variants_SNP_indel")

# Clone the header formatting (fill + border) from the last existing header
# cell (G1, "Comments") onto the new header cell (H1) so it reuses the same
# cell style instead of minting a new one.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# New column header text (this also grows the shared-string table).
$ws.Range("H1").Value = "Reference/s"

# Give the new column the same authored width as the rest of the header row.
$ws.Columns.Item(8).ColumnWidth = 28.5

# Document the new field the same way the other header cells are documented.
$commentText = "Derek:" + "`n" + "Publication where the mutation/resistance determinant was first identified. Free form field"
$ws.Range("H1").AddComment($commentText) | Out-Null

# Restore the selection to where the author last left it.
$ws.Range("H18").Select() | Out-Null
